$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.754.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.600.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.622.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "

$ws.Range("E11").Value = "  +2.65%  "

$ws.Range("E13").Value = "  +7.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.063.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.760.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.28%  "

$ws.Range("E17").Value = "  +2.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.615.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.85%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.81%  "

$ws.Range("E22").Value = "  +12.42%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.524"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0791"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.77%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +4.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.79%  "

$ws.Range("E35").Value = "  +5.50%  "

$ws.Range("E36").Value = "  +9.48%  "

$ws.Range("E37").Value = "  +3.52%  "

$ws.Range("E38").Value = "  +7.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.849"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.84%  "

$ws.Range("E41").Value = "  +3.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0546"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.10%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.36%  "

